$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.078.26"
$ws.Range("E2").Value = "  +1.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.996.63"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.50"
$ws.Range("E5").Value = "  +15.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.63"
$ws.Range("E6").Value = "  +10.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.688"
$ws.Range("E7").Value = "  -1.34%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.763"
$ws.Range("E9").Value = "  +1.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  -2.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.36"
$ws.Range("E11").Value = "  +8.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000319"
$ws.Range("E12").Value = "  -2.05%  "

$ws.Range("E13").Value = "  +4.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.632.62"
$ws.Range("E14").Value = "  -0.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.997.24"
$ws.Range("E15").Value = "  -0.76%  "

$ws.Range("E16").Value = "  +6.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.38"
$ws.Range("E17").Value = "  +1.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.72"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.931.21"
$ws.Range("E20").Value = "  +1.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "443.31"
$ws.Range("E21").Value = "  +2.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.96"
$ws.Range("E22").Value = "  +17.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.69"
$ws.Range("E23").Value = "  -1.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.41"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.50"
$ws.Range("E25").Value = "  -1.66%  "

$ws.Range("E26").Value = "  -5.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.36"
$ws.Range("E27").Value = "  -0.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.60"
$ws.Range("E28").Value = "  -1.83%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  -2.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.88"
$ws.Range("E31").Value = "  -4.33%  "

$ws.Range("E32").Value = "  +2.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.132"
$ws.Range("E33").Value = "  -3.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "49.24"
$ws.Range("E34").Value = "  -1.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "71.73"
$ws.Range("E35").Value = "  +5.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "641.64"
$ws.Range("E36").Value = "  -5.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0917"
$ws.Range("E37").Value = "  +11.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.437"
$ws.Range("E38").Value = "  -5.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.51"
$ws.Range("E39").Value = "  +3.52%  "

$ws.Range("E40").Value = "  -0.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.03"
$ws.Range("E42").Value = "  -0.85%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.31"
$ws.Range("E43").Value = "  -3.25%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("E46").Value = "  -0.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.98"
$ws.Range("E47").Value = "  +37.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.64"
$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.40"
$ws.Range("E49").Value = "  +0.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.910.43"
$ws.Range("E50").Value = "  +1.72%  "

$ws.Range("E51").Value = "  +0.23%  "
